$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure Price column cells keep their original Text format so numeric-looking
# values (e.g. "228.54") are not reinterpreted as numbers by Excel.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "38.726.77"
$ws.Range("E2").Value = "  +1.06%  "
$ws.Range("D3").Value = "2.101.63"
$ws.Range("E3").Value = "  +0.30%  "
$ws.Range("E4").Value = "  +0.01%  "
$ws.Range("D5").Value = "228.54"
$ws.Range("E5").Value = "  +0.02%  "
$ws.Range("E6").Value = "  +0.48%  "
$ws.Range("D7").Value = "62.14"
$ws.Range("E7").Value = "  +1.74%  "
$ws.Range("E8").Value = "  -0.04%  "
$ws.Range("E9").Value = "  +2.30%  "
$ws.Range("E10").Value = "  -0.30%  "
$ws.Range("E11").Value = "  -0.35%  "
$ws.Range("D12").Value = "15.82"
$ws.Range("E12").Value = "  +7.04%  "
$ws.Range("D13").Value = "2.413.34"
$ws.Range("E13").Value = "  +0.34%  "
$ws.Range("D14").Value = "22.06"
$ws.Range("E14").Value = "  -1.21%  "
$ws.Range("E15").Value = "  +3.76%  "
$ws.Range("D16").Value = "5.54"
$ws.Range("E16").Value = "  +1.71%  "
$ws.Range("D17").Value = "2.104.57"
$ws.Range("E17").Value = "  +0.40%  "
$ws.Range("D18").Value = "38.743.26"
$ws.Range("E18").Value = "  +1.24%  "
$ws.Range("D19").Value = "71.87"
$ws.Range("E20").Value = "  +1.09%  "
$ws.Range("E21").Value = "  +0.90%  "
$ws.Range("D22").Value = "227.61"
$ws.Range("E22").Value = "  +1.11%  "
$ws.Range("D24").Value = "2.35"
$ws.Range("E24").Value = "  -3.57%  "
$ws.Range("D25").Value = "2.32"
$ws.Range("E25").Value = "  +0.36%  "
$ws.Range("D26").Value = "172.35"
$ws.Range("E26").Value = "  +1.54%  "
$ws.Range("D27").Value = "9.58"
$ws.Range("E27").Value = "  +1.83%  "
$ws.Range("E28").Value = "  +6.24%  "
$ws.Range("E29").Value = "  +4.15%  "
$ws.Range("E30").Value = "  +1.74%  "
$ws.Range("D31").Value = "2.51"
$ws.Range("E31").Value = "  +6.64%  "
$ws.Range("E32").Value = "  +1.39%  "
$ws.Range("D33").Value = "4.54"
$ws.Range("E33").Value = "  +2.40%  "
$ws.Range("E34").Value = "  -0.13%  "
$ws.Range("E35").Value = "  +2.64%  "
$ws.Range("D36").Value = "6.90"
$ws.Range("E36").Value = "  +7.31%  "
$ws.Range("E37").Value = "  +1.41%  "
$ws.Range("D38").Value = "3.58"
$ws.Range("E38").Value = "  +2.23%  "
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").Value = "18.13"
$ws.Range("E40").Value = "  -0.65%  "
$ws.Range("D41").Value = "102.78"
$ws.Range("E41").Value = "  +2.99%  "
$ws.Range("E42").Value = "  +4.37%  "
$ws.Range("D43").Value = "1.534.02"
$ws.Range("E43").Value = "  -0.51%  "
$ws.Range("D44").Value = "7.85"
$ws.Range("E44").Value = "  +4.65%  "
$ws.Range("E45").Value = "  -1.19%  "
$ws.Range("E46").Value = "  +3.34%  "
$ws.Range("D47").Value = "0.0911"
$ws.Range("E47").Value = "  +0.27%  "
$ws.Range("E48").Value = "  +0.02%  "
$ws.Range("E49").Value = "  +1.58%  "
$ws.Range("E50").Value = "  -0.31%  "
$ws.Range("D51").Value = "2.294.32"
$ws.Range("E51").Value = "  +0.14%  "
